$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 41 (pushes the existing rows 41..104 down to 42..105)
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with its values
$ws.Cells.Item(41, 1).Value = 11
$ws.Cells.Item(41, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(41, 3).Value = "Bíobío"
$ws.Cells.Item(41, 4).Value = 44477
$ws.Cells.Item(41, 5).Value = 8
$ws.Cells.Item(41, 6).Value = "Fruta"
$ws.Cells.Item(41, 7).Value = 100108
$ws.Cells.Item(41, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(41, 9).Value = 100108005
$ws.Cells.Item(41, 10).Value = "Piña"
$ws.Cells.Item(41, 11).Value = "Caramelo"
$ws.Cells.Item(41, 12).Value = "Segunda"
$ws.Cells.Item(41, 13).Value = 100
$ws.Cells.Item(41, 14).Value = 20000
$ws.Cells.Item(41, 15).Value = 21000
$ws.Cells.Item(41, 16).Value = 20500
$ws.Cells.Item(41, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(41, 18).Value = "Ecuador"
$ws.Cells.Item(41, 19).Value = 1464
$ws.Cells.Item(41, 20).Value = 14

# Keep the date column formatted the same way as the rest of column D
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
